# Commit "Add files via upload": BOM.xlsx re-uploaded with the header
# label in A1 corrected from "class_name" to "Class_Name".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")
$ws.Range("A1").Value = "Class_Name"
